$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the "(c) 2020 ..." footer
# paragraph by their text content (more robust than hard-coded indices).
# The blank paragraph immediately preceding "Ver no Jupiter" is removed
# together with both text paragraphs, while the blank paragraph that
# follows the footer (and the page-break paragraph after it) are kept.
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "*Ver no Jupiter*") {
        $startIdx = $i - 1
    }
    if ($txt -like "*Contact: luizeleno@usp.br*") {
        $endIdx = $i
    }
}

if ($startIdx -ge 1 -and $endIdx -ge $startIdx) {
    $startPara = $d.Paragraphs.Item($startIdx)
    $endPara = $d.Paragraphs.Item($endIdx)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
